# "fix template, label satuan dan tampilan"
#
# The single wide header row (A1:CF1 -- 84 columns: mp1-6, mu1-6, fr1-24,
# amp1-24, gh1-24) is split into four shorter header blocks, each
# followed by an empty bordered data row:
#   row 1  (A:L)  mp1-6 / mu1-6     -> row 2  (A:L)  blank data (quote-prefixed)
#   row 4  (A:X)  fr1-24            -> row 5  (A:X)  blank data
#   row 7  (A:X)  amp1-24           -> row 8  (A:X)  blank data
#   row 10 (A:X)  gh1-24            -> row 11 (A:X)  blank data
# A "Catatan" (notes) block is added at A14:A16 documenting the
# apostrophe-before-zero convention used in row 2. The view is re-zoomed
# to 125% and the selection moved to I19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Col-Letter([int]$i) {
    $s = ""
    while ($i -gt 0) {
        $m = ($i - 1) % 26
        $s = [char](65 + $m) + $s
        $i = [int](($i - $m - 1) / 26)
    }
    return $s
}

# ---------------------------------------------------------------------
# 0. Read the three 24-wide label blocks off the original wide row 1
#    (columns M:AK = fr1-24, AL:BH = amp1-24, BI:CF = gh1-24) before we
#    touch anything.
# ---------------------------------------------------------------------
$fr = @()
for ($c = 13; $c -le 36; $c++) { $fr += $ws.Cells.Item(1, $c).Text }

$amp = @()
for ($c = 37; $c -le 60; $c++) { $amp += $ws.Cells.Item(1, $c).Text }

$gh = @()
for ($c = 61; $c -le 84; $c++) { $gh += $ws.Cells.Item(1, $c).Text }


# ---------------------------------------------------------------------
# 1. Wipe out the old wide tail (M1:CF2) -- the mp/mu block at A1:L1 (and
#    its A2:L2 data row) stays put, everything east of column L goes away
#    so we can rebuild fr/amp/gh further down the sheet instead.
# ---------------------------------------------------------------------
$ws.Range("M1:CF2").Clear() | Out-Null

# ---------------------------------------------------------------------
# 2. Helper to stamp a 24-wide header row (A:X) with a label array using
#    the same visual style as the existing A1:L1 header (bold + border),
#    and an empty bordered "data" row beneath it (style copied from a
#    caller-supplied source cell/range so rows 5 and 2 can share the
#    quote-prefix-capable look while rows 8/11 get the plain one).
# ---------------------------------------------------------------------
function Write-HeaderBlock([int]$headerRow, [int]$dataRow, $labels, [string]$dataStyleSource) {
    $lastCol = Col-Letter $labels.Length
    # Header row: copy the bold+bordered look of A1 across A:<lastCol>, then
    # fill in the text.
    $ws.Range("A1").Copy() | Out-Null
    $ws.Range("A${headerRow}:${lastCol}${headerRow}").PasteSpecial(-4122) | Out-Null
    for ($i = 0; $i -lt $labels.Length; $i++) {
        $col = Col-Letter ($i + 1)
        $ws.Cells.Item($headerRow, $i + 1).Value = $labels[$i]
    }

    # Blank bordered data row underneath.
    $ws.Range($dataStyleSource).Copy() | Out-Null
    $ws.Range("A${dataRow}:${lastCol}${dataRow}").PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# 3. Stash the plain bordered look (no quote-prefix) on a scratch cell
#    *before* row 2 is converted to quote-prefixed text below, so rows 8
#    and 11 can still be stamped with the original plain style.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 4. Turn row 2 (A:L) into the quote-prefixed blank data row: an
#    apostrophe-only entry stores an empty, quote-prefixed string so a
#    literal "0" typed later displays as text instead of a number.
# ---------------------------------------------------------------------
$ws.Range("A2:L2").Value = "'"

Write-HeaderBlock 4 5 $fr "A2"
Write-HeaderBlock 7 8 $amp "Z1"
Write-HeaderBlock 10 11 $gh "Z1"

$ws.Range("Z1").Clear() | Out-Null

# ---------------------------------------------------------------------
# 5. "Catatan" (notes) block: a bold, slightly larger title followed by
#    two italic explanatory lines (Indonesian + English) describing the
#    apostrophe-before-zero convention used above.
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Catatan"
$ws.Range("A14").Font.Bold = $true
$ws.Range("A14").Font.Size = 12

$ws.Range("A15").Value = "Jika nilai adalah data adalah 0 maka beri tanda petik 1 sebelum angka 0"
$ws.Range("A15").Font.Italic = $true

$ws.Range("A16").Value = "If the data value is 0 then put a quotation mark 1 before the number 0"
$ws.Range("A16").Font.Italic = $true

# ---------------------------------------------------------------------
# 6. View tweaks: zoom to 125% and move the active selection to I19.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 125
$ws.Range("I19").Select() | Out-Null

